# "Generate Report for Archive"
# The localization status moved from "Ready for handoff" to "In Translation"
# on every sheet that tracks it:
#   - Overview sheet: columns E (zh-cn) and F (de-de), row 2
#   - zh-cn sheet:     column C (Status), row 2
#   - de-de sheet:     column C (Status), row 2
# After the text shrank, the Status column was re-autofit to the new
# (narrower) text width.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
# Width (in "characters") that the Status columns settle on once they are
# auto-fit to the shorter "In Translation" label.
$newStatusColumnWidth = 12.5

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColumnWidth

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColumnWidth

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColumnWidth
